$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price cells are forced to remain text (matching the original inlineStr cells)
# by temporarily applying a text NumberFormat, then the style is reset to Normal
# so no residual formatting is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.619.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.904.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.62%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.506'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.902.61'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.49%  '

$ws.Range("E11").Value = '  -3.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.37%  '

$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.43%  '

$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.385.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.570.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.903.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '426.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.91%  '

$ws.Range("E22").Value = '  -2.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.13'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.82%  '

$ws.Range("E25").Value = '  +1.37%  '

$ws.Range("E26").Value = '  -1.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.20%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.63%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.37%  '

$ws.Range("E32").Value = '  -2.97%  '

$ws.Range("E33").Value = '  -3.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0836'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.63%  '

$ws.Range("E36").Value = '  -2.11%  '

$ws.Range("E37").Value = '  -2.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("E39").Value = '  +0.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.39'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.54%  '

$ws.Range("E41").Value = '  -3.86%  '

$ws.Range("E42").Value = '  -1.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.291'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0347'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '371.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.653.92'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  -0.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.57%  '

$ws.Range("E51").Value = '  -1.10%  '
